$d = $word.ActiveDocument

# 1. Justify ("both") every paragraph in the body except the document Title.
$count = $d.Paragraphs.Count
for ($i = 1; $i -le $count; $i++) {
    $para = $d.Paragraphs($i)
    $styleName = $para.Range.ParagraphFormat.Style.NameLocal
    if ($styleName -ne "Title") {
        $para.Range.ParagraphFormat.Alignment = 3
    }
}

# 2. Remove the old "_GoBack" bookmark that sits between
#    "...Statistical G" and "enetics, or similar." in the Education sentence.
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks("_GoBack").Delete()
}

# 3. Re-create the "_GoBack" bookmark at the new location: splitting the
#    final paragraph's run right after "...characteristic prot" and before
#    "ected by law. ".
$rng = $d.Content
$rng.Find.Execute("characteristic prot", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$collapsedRng = $rng.Duplicate
$collapsedRng.Collapse(0)
$d.Bookmarks.Add("_GoBack", $collapsedRng)

Write-Output "done"
